$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The export window slid forward by 2 days: the first two data rows
# (2025-11-21 and 2025-11-22) fall out of range, and every remaining
# row's data shifts up by two rows (row N now holds what used to be at
# row N+2), shrinking the table from A1:D89 to A1:D87.
$ws.Rows("2:3").Delete()
